# "Adding Login Module and Word Add In"
# Updates the Login sheet credentials used by the Ranorex/automation suite
# (firm id, username, server name) and moves the active selection, to
# reflect the new QA environment the login module now points at.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")

# Set values in first-use order: firmid, QA User1, newautomation, Amicus User
$ws.Range("A2").Value2 = "firmid"
$ws.Range("A3").Value2 = "firmid"
$ws.Range("B3").Value2 = "QA User1"
$ws.Range("D2").Value2 = "newautomation"
$ws.Range("D3").Value2 = "newautomation"
$ws.Range("B2").Value2 = "Amicus User"

# Move the sheet's active cell/selection from B4 to C6
$null = $ws.Range("C6").Select()
